# Auto-generated: bulk value refresh across Sheets (scheduled runner update)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1202.49
$ws.Range("I15").Value = 1202.49
$ws.Range("K15").Value = 3607.47
$ws.Range("M15").Value = -3438.47
$ws.Range("H17").Value = 1003.4
$ws.Range("J17").Value = 1003.4
$ws.Range("L17").Value = 3010.2
$ws.Range("N17").Value = -3346.2
$ws.Range("H98").Value = 8298.647000000001
$ws.Range("I98").Value = 9648.5
$ws.Range("J98").Value = 1999.3334
$ws.Range("K98").Value = 9648.5
$ws.Range("L98").Value = 1999.3334
$ws.Range("M98").Value = -8150.5
$ws.Range("N98").Value = -4995.3334
$ws.Range("H100").Value = 2084.2856
$ws.Range("I100").Value = 1719.8
$ws.Range("J100").Value = 2995.5
$ws.Range("K100").Value = 1719.8
$ws.Range("L100").Value = 2995.5
$ws.Range("M100").Value = -1178.8
$ws.Range("N100").Value = -4077.5
$ws.Range("H122").Value = 8298.647000000001
$ws.Range("I122").Value = 9648.5
$ws.Range("J122").Value = 1999.3334
$ws.Range("K122").Value = 28945.5
$ws.Range("L122").Value = 5998.0002
$ws.Range("M122").Value = -26495.5
$ws.Range("N122").Value = -10898.0002
$ws.Range("H123").Value = 28999.5
$ws.Range("J123").Value = 28999.5
$ws.Range("L123").Value = 28999.5
$ws.Range("N123").Value = -38799.5
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H138").Value = 2796.5115
$ws.Range("I138").Value = 2718.2
$ws.Range("J138").Value = 2806.5513
$ws.Range("K138").Value = 8154.599999999999
$ws.Range("L138").Value = 8419.653900000001
$ws.Range("M138").Value = -3014.599999999999
$ws.Range("N138").Value = -18699.6539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13263.82
$ws.Range("I32").Value = 9147.132
$ws.Range("J32").Value = 17906.043
$ws.Range("K32").Value = 9147.132
$ws.Range("L32").Value = 17906.043
$ws.Range("M32").Value = -8860.132
$ws.Range("N32").Value = -18480.043
$ws.Range("H102").Value = 9807053
$ws.Range("I102").Value = 9807053
$ws.Range("K102").Value = 9807053
$ws.Range("M102").Value = -9805431
$ws.Range("H132").Value = 2965.1292
$ws.Range("I132").Value = 2522.652
$ws.Range("K132").Value = 7567.956
$ws.Range("M132").Value = -5037.956

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H105").Value = 500001000
$ws.Range("I105").Value = 500001000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 500001000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -499999253
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 52632732
$ws.Range("I16").Value = 58824640
$ws.Range("K16").Value = 58824640
$ws.Range("M16").Value = -58824353
$ws.Range("H31").Value = 1694.746
$ws.Range("I31").Value = 1621.2075
$ws.Range("J31").Value = 2084.5
$ws.Range("K31").Value = 1621.2075
$ws.Range("L31").Value = 2084.5
$ws.Range("M31").Value = -1326.2075
$ws.Range("N31").Value = -2674.5
$ws.Range("H34").Value = 1694.746
$ws.Range("I34").Value = 1621.2075
$ws.Range("J34").Value = 2084.5
$ws.Range("K34").Value = 1621.2075
$ws.Range("L34").Value = 2084.5
$ws.Range("M34").Value = -1419.2075
$ws.Range("N34").Value = -2488.5
$ws.Range("H99").Value = 2025940.5
$ws.Range("I99").Value = 6580432
$ws.Range("J99").Value = 1722.1111
$ws.Range("K99").Value = 6580432
$ws.Range("L99").Value = 1722.1111
$ws.Range("M99").Value = -6578934
$ws.Range("N99").Value = -4718.1111
$ws.Range("H113").Value = 52632732
$ws.Range("I113").Value = 58824640
$ws.Range("K113").Value = 58824640
$ws.Range("M113").Value = -58822470
$ws.Range("H126").Value = 2025940.5
$ws.Range("I126").Value = 6580432
$ws.Range("J126").Value = 1722.1111
$ws.Range("K126").Value = 19741296
$ws.Range("L126").Value = 5166.3333
$ws.Range("M126").Value = -19738826
$ws.Range("N126").Value = -10106.3333
$ws.Range("H141").Value = 245226.89
$ws.Range("J141").Value = 245226.89
$ws.Range("L141").Value = 245226.89
$ws.Range("N141").Value = -255586.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 155212.55
$ws.Range("I11").Value = 182428.3
$ws.Range("K11").Value = 547284.8999999999
$ws.Range("M11").Value = -547144.8999999999
$ws.Range("H56").Value = 6437.6665
$ws.Range("I56").Value = 6437.6665
$ws.Range("K56").Value = 6437.6665
$ws.Range("M56").Value = -5907.6665
$ws.Range("H68").Value = 1458.5
$ws.Range("I68").Value = 1300.5
$ws.Range("J68").Value = 1537.5
$ws.Range("K68").Value = 3901.5
$ws.Range("L68").Value = 4612.5
$ws.Range("M68").Value = -3090.5
$ws.Range("N68").Value = -6234.5
$ws.Range("H69").Value = 1666.2
$ws.Range("I69").Value = 1462
$ws.Range("K69").Value = 4386
$ws.Range("M69").Value = -3575
$ws.Range("H71").Value = 1458.5
$ws.Range("I71").Value = 1300.5
$ws.Range("J71").Value = 1537.5
$ws.Range("K71").Value = 11704.5
$ws.Range("L71").Value = 13837.5
$ws.Range("M71").Value = -7648.5
$ws.Range("N71").Value = -21949.5
$ws.Range("H72").Value = 1666.2
$ws.Range("I72").Value = 1462
$ws.Range("K72").Value = 13158
$ws.Range("M72").Value = -9102
$ws.Range("H80").Value = 2985
$ws.Range("I80").Value = 2688.3333
$ws.Range("J80").Value = 3133.3333
$ws.Range("K80").Value = 8064.999899999999
$ws.Range("L80").Value = 9399.999899999999
$ws.Range("M80").Value = -7128.999899999999
$ws.Range("N80").Value = -11271.9999
$ws.Range("H83").Value = 2985
$ws.Range("I83").Value = 2688.3333
$ws.Range("J83").Value = 3133.3333
$ws.Range("K83").Value = 24194.9997
$ws.Range("L83").Value = 28199.9997
$ws.Range("M83").Value = -19514.9997
$ws.Range("N83").Value = -37559.9997
$ws.Range("H122").Value = 1026.9166
$ws.Range("J122").Value = 1101.5
$ws.Range("L122").Value = 9913.5
$ws.Range("N122").Value = -14813.5
$ws.Range("H131").Value = 26356276
$ws.Range("I131").Value = 166667100
$ws.Range("J131").Value = 47994.03
$ws.Range("K131").Value = 500001300
$ws.Range("L131").Value = 143982.09
$ws.Range("M131").Value = -499996260
$ws.Range("N131").Value = -154062.09
$ws.Range("H132").Value = 1360.3125
$ws.Range("I132").Value = 1240.4445
$ws.Range("J132").Value = 1514.4286
$ws.Range("K132").Value = 11164.0005
$ws.Range("L132").Value = 13629.8574
$ws.Range("M132").Value = -8634.0005
$ws.Range("N132").Value = -18689.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3475
$ws.Range("I80").Value = 1840
$ws.Range("K80").Value = 1840
$ws.Range("M80").Value = -842
$ws.Range("H83").Value = 3475
$ws.Range("I83").Value = 1840
$ws.Range("K83").Value = 9200
$ws.Range("M83").Value = -4208
$ws.Range("H113").Value = 1820
$ws.Range("I113").Value = 1980
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1980
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 190
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 7103.769
$ws.Range("I132").Value = 7904.3
$ws.Range("J132").Value = 4435.3335
$ws.Range("K132").Value = 23712.9
$ws.Range("L132").Value = 13306.0005
$ws.Range("M132").Value = -21182.9
$ws.Range("N132").Value = -18366.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1718.091
$ws.Range("I61").Value = 1312.375
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 1312.375
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -1110.375
$ws.Range("N61").Value = -3204
$ws.Range("H82").Value = 1978.2727
$ws.Range("I82").Value = 2086.111
$ws.Range("J82").Value = 1493
$ws.Range("K82").Value = 2086.111
$ws.Range("L82").Value = 1493
$ws.Range("M82").Value = -1725.111
$ws.Range("N82").Value = -2215
$ws.Range("H85").Value = 1978.2727
$ws.Range("I85").Value = 2086.111
$ws.Range("J85").Value = 1493
$ws.Range("K85").Value = 2086.111
$ws.Range("L85").Value = 1493
$ws.Range("M85").Value = -838.1109999999999
$ws.Range("N85").Value = -3989
$ws.Range("H100").Value = 1399.7142
$ws.Range("I100").Value = 1159.6
$ws.Range("K100").Value = 1159.6
$ws.Range("M100").Value = -618.5999999999999
$ws.Range("H113").Value = 1718.091
$ws.Range("I113").Value = 1312.375
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 1312.375
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = 857.625
$ws.Range("N113").Value = -7140
$ws.Range("H132").Value = 3868.2
$ws.Range("I132").Value = 4562.3335
$ws.Range("J132").Value = 3570.7144
$ws.Range("K132").Value = 13687.0005
$ws.Range("L132").Value = 10712.1432
$ws.Range("M132").Value = -11157.0005
$ws.Range("N132").Value = -15772.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 19233348
$ws.Range("I122").Value = 22729856
$ws.Range("K122").Value = 68189568
$ws.Range("M122").Value = -68187118
$ws.Range("H136").Value = 1733.3334
$ws.Range("I136").Value = 1320
$ws.Range("J136").Value = 2028.5714
$ws.Range("K136").Value = 3960
$ws.Range("L136").Value = 6085.7142
$ws.Range("M136").Value = -1410
$ws.Range("N136").Value = -11185.7142
$ws.Range("H141").Value = 97450
$ws.Range("J141").Value = 97450
$ws.Range("L141").Value = 97450
$ws.Range("N141").Value = -107810

